$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.374.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.08%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.787.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.34%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'335.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.73%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.10%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.3786"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.33%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3417"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.04%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'48.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.10%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'1.200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.17%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07458"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.07%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.18%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'21.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +8.83%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.459"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.70%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'1.787.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.75%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'7.014"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.67%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.00001091"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.42%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.06627"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.52%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'84.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.69%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.16%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'17.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +4.52%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.461"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.73%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'27.332.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.05%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'12.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.32%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.450"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.04%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.545"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.84%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'1.481"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.21%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'21.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +9.26%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'150.19"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").Value = "'1.987.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.83%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'132.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.17%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'4.068"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.78%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'6.100"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.74%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.08677"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.77%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'13.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.99%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.671"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.60%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.6864"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +10.36%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'5.413"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.55%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.06331"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.75%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'8.802"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.26%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'Algorand"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.2189"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.34%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'VeChain"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.02336"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.15%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.274"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.21%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'14.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.13%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.002"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.23%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.6415"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +6.32%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'3.847"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.65%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'2.108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.74%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'129.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.15%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.07182"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.04%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'78.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.32%  "
$ws.Range("E51").Style = "Normal"
